$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "563.45") must be
# forced to Text format first, otherwise Excel auto-converts them to a
# floating point number (losing the exact decimal string / introducing
# floating point noise).
$ws.Range('D2').Value = '62.189.20'
$ws.Range('E2').Value = '  +1.39%  '
$ws.Range('D3').Value = '2.421.90'
$ws.Range('E3').Value = '  +1.74%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.45'
$ws.Range('E5').Value = '  +2.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.16'
$ws.Range('E6').Value = '  +3.58%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +1.74%  '
$ws.Range('D9').Value = '2.419.43'
$ws.Range('E9').Value = '  +1.63%  '
$ws.Range('E10').Value = '  +1.21%  '
$ws.Range('E11').Value = '  -1.59%  '
$ws.Range('E12').Value = '  +1.44%  '
$ws.Range('E13').Value = '  +1.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.11'
$ws.Range('E14').Value = '  +3.60%  '
$ws.Range('E15').Value = '  +5.42%  '
$ws.Range('D16').Value = '2.853.76'
$ws.Range('D17').Value = '61.949.59'
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('D18').Value = '2.416.80'
$ws.Range('E18').Value = '  +1.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.21'
$ws.Range('E19').Value = '  +2.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '324.77'
$ws.Range('E20').Value = '  +1.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.18'
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.76'
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.44'
$ws.Range('E24').Value = '  +1.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.72'
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.04'
$ws.Range('E26').Value = '  +5.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '588.34'
$ws.Range('E27').Value = '  +14.58%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.527.00'
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0943'
$ws.Range('E29').Value = '  +5.10%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.994'
$ws.Range('E30').Value = '  -0.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.26'
$ws.Range('E31').Value = '  +0.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.45'
$ws.Range('E32').Value = '  +5.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.149'
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('E34').Value = '  +2.72%  '
$ws.Range('E35').Value = '  +1.95%  '
$ws.Range('E36').Value = '  +4.67%  '
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.77'
$ws.Range('E38').Value = '  +1.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '154.42'
$ws.Range('E39').Value = '  +5.28%  '
$ws.Range('E40').Value = '  +1.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.69'
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.83'
$ws.Range('E42').Value = '  -5.01%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.37'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '150.62'
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('E46').Value = '  +1.43%  '
$ws.Range('E47').Value = '  +3.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.36'
$ws.Range('E48').Value = '  +4.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.592'
$ws.Range('E49').Value = '  +2.25%  '
$ws.Range('E50').Value = '  +1.97%  '
$ws.Range('E51').Value = '  +2.31%  '
